# Apply the changes described by the diff:
#  - Update K16, K19, K22, K23, K27 on "Customer Quote" sheet from 1.0565 to 1
#  - Update the active selection from A27 to G9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

$ws.Range("K16").Value = 1
$ws.Range("K19").Value = 1
$ws.Range("K22").Value = 1
$ws.Range("K23").Value = 1
$ws.Range("K27").Value = 1

$ws.Activate()
$ws.Range("G9").Select()
